# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# Updates numeric price/profit columns (H-N) for specific rows across multiple sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 914
$ws.Range("I41").Value = 759.6
$ws.Range("J41").Value = 1300
$ws.Range("K41").Value = 759.6
$ws.Range("L41").Value = 1300
$ws.Range("M41").Value = -319.6
$ws.Range("N41").Value = -2180
$ws.Range("H86").Value = 71462140
$ws.Range("I86").Value = 4452.7
$ws.Range("K86").Value = 4452.7
$ws.Range("M86").Value = -3329.7
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H89").Value = 71462140
$ws.Range("I89").Value = 4452.7
$ws.Range("K89").Value = 22263.5
$ws.Range("M89").Value = -16647.5
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H113").Value = 10931.471
$ws.Range("I113").Value = 13292.625
$ws.Range("K113").Value = 13292.625
$ws.Range("M113").Value = -10038.625
$ws.Range("H137").Value = 629424.5
$ws.Range("I137").Value = 837566.25
$ws.Range("J137").Value = 4999.25
$ws.Range("K137").Value = 2512698.75
$ws.Range("L137").Value = 14997.75
$ws.Range("M137").Value = -2510148.75
$ws.Range("N137").Value = -20097.75
$ws.Range("H138").Value = 417889.5
$ws.Range("J138").Value = 7033.28
$ws.Range("L138").Value = 21099.84
$ws.Range("N138").Value = -31379.84

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24790.105
$ws.Range("I32").Value = 18883.059
$ws.Range("J32").Value = 75000
$ws.Range("K32").Value = 18883.059
$ws.Range("L32").Value = 75000
$ws.Range("M32").Value = -18596.059
$ws.Range("N32").Value = -75574

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5206.6924
$ws.Range("I20").Value = 2837.5
$ws.Range("K20").Value = 2837.5
$ws.Range("M20").Value = -2590.5
$ws.Range("H58").Value = 53044
$ws.Range("J58").Value = 58433.168
$ws.Range("L58").Value = 58433.168
$ws.Range("N58").Value = -59021.168

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 450
$ws.Range("I3").Value = 450
$ws.Range("K3").Value = 450
$ws.Range("M3").Value = -337
$ws.Range("H7").Value = 974.6
$ws.Range("I7").Value = 624.8333
$ws.Range("K7").Value = 624.8333
$ws.Range("M7").Value = -511.8333
$ws.Range("H31").Value = 4661.5
$ws.Range("I31").Value = 4212.5713
$ws.Range("J31").Value = 5290
$ws.Range("K31").Value = 4212.5713
$ws.Range("L31").Value = 5290
$ws.Range("M31").Value = -3917.5713
$ws.Range("N31").Value = -5880
$ws.Range("H34").Value = 4661.5
$ws.Range("I34").Value = 4212.5713
$ws.Range("J34").Value = 5290
$ws.Range("K34").Value = 4212.5713
$ws.Range("L34").Value = 5290
$ws.Range("M34").Value = -4010.5713
$ws.Range("N34").Value = -5694
$ws.Range("H58").Value = 3144.9443
$ws.Range("I58").Value = 2145.5386
$ws.Range("K58").Value = 2145.5386
$ws.Range("M58").Value = -1942.5386
$ws.Range("H132").Value = 6933.579
$ws.Range("I132").Value = 7783.6875
$ws.Range("K132").Value = 23351.0625
$ws.Range("M132").Value = -20821.0625
$ws.Range("H136").Value = 3144.9443
$ws.Range("I136").Value = 2145.5386
$ws.Range("K136").Value = 6436.6158
$ws.Range("M136").Value = -3886.6158
$ws.Range("H141").Value = 568932.9399999999
$ws.Range("J141").Value = 598241.25
$ws.Range("L141").Value = 598241.25
$ws.Range("N141").Value = -608601.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 45752560
$ws.Range("I4").Value = 45408250
$ws.Range("J4").Value = 46747240
$ws.Range("K4").Value = 136224750
$ws.Range("L4").Value = 140241720
$ws.Range("M4").Value = -136224638
$ws.Range("N4").Value = -140241944
$ws.Range("H11").Value = 601705.25
$ws.Range("J11").Value = 167150.33
$ws.Range("L11").Value = 501450.99
$ws.Range("N11").Value = -501730.99
$ws.Range("H80").Value = 47076480
$ws.Range("I80").Value = 250002000
$ws.Range("J80").Value = 10180930
$ws.Range("K80").Value = 750006000
$ws.Range("L80").Value = 30542790
$ws.Range("M80").Value = -750005064
$ws.Range("N80").Value = -30544662
$ws.Range("H83").Value = 47076480
$ws.Range("I83").Value = 250002000
$ws.Range("J83").Value = 10180930
$ws.Range("K83").Value = 2250018000
$ws.Range("L83").Value = 91628370
$ws.Range("M83").Value = -2250013320
$ws.Range("N83").Value = -91637730
$ws.Range("H122").Value = 6275.276
$ws.Range("J122").Value = 8129.6665
$ws.Range("L122").Value = 73166.9985
$ws.Range("N122").Value = -78066.9985
$ws.Range("H131").Value = 5773.409
$ws.Range("I131").Value = 7736.5454
$ws.Range("J131").Value = 3810.2727
$ws.Range("K131").Value = 23209.6362
$ws.Range("L131").Value = 11430.8181
$ws.Range("M131").Value = -18169.6362
$ws.Range("N131").Value = -21510.8181

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10699
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 10699
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 10699
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -11239
$ws.Range("H73").Value = 10699
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 10699
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 10699
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -12571
$ws.Range("I80").Value = 38749.75
$ws.Range("J80").Value = 10000
$ws.Range("K80").Value = 38749.75
$ws.Range("L80").Value = 10000
$ws.Range("M80").Value = -37751.75
$ws.Range("N80").Value = -11996
$ws.Range("I83").Value = 38749.75
$ws.Range("J83").Value = 10000
$ws.Range("K83").Value = 193748.75
$ws.Range("L83").Value = 50000
$ws.Range("M83").Value = -188756.75
$ws.Range("N83").Value = -59984
$ws.Range("H126").Value = 11346.125
$ws.Range("I126").Value = 16288.889
$ws.Range("K126").Value = 48866.667
$ws.Range("M126").Value = -46396.667
$ws.Range("H132").Value = 3347
$ws.Range("I132").Value = 3502.5881
$ws.Range("J132").Value = 1760
$ws.Range("K132").Value = 10507.7643
$ws.Range("L132").Value = 5280
$ws.Range("M132").Value = -7977.764299999999
$ws.Range("N132").Value = -10340

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 66000
$ws.Range("J46").Value = 68000
$ws.Range("L46").Value = 68000
$ws.Range("N46").Value = -68462
$ws.Range("H100").Value = 25742
$ws.Range("I100").Value = 5841.6665
$ws.Range("J100").Value = 51328.145
$ws.Range("K100").Value = 11683.333
$ws.Range("L100").Value = 102656.29
$ws.Range("M100").Value = -11142.333
$ws.Range("N100").Value = -103738.29
$ws.Range("H134").Value = 66000
$ws.Range("J134").Value = 68000
$ws.Range("L134").Value = 204000
$ws.Range("N134").Value = -209070
